# Delete old non-pharma patients and add excluded cases.
#
# Target sheet: "NonPharma Interventions (NPI)" (3rd sheet / sheet3.xml).
# A new column is inserted at D ("Created"), pushing the old D..M columns to
# E..N. The new column is populated with "Created"/"X" markers on the first
# few (header/excluded) rows, and left blank everywhere else. The sheet also
# becomes the active sheet/tab, with the active cell parked on D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NonPharma Interventions (NPI)")

# Insert a new blank column before the existing column D ("PatientID" data
# column), shifting D:M to E:N.
$ws.Columns.Item(4).Insert()

# Narrower, manually-set (non bestFit) width for the new column, matching
# the author's resize of the "Created" column.
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334

# New column header.
$ws.Range("D1").Value = "Created"

# Mark the "X" rows for the new column.
$ws.Range("D2").Value = "X"
$ws.Range("D3").Value = "X"
$ws.Range("D5").Value = "X"
$ws.Range("D6").Value = "X"

# All remaining rows keep the new column fully empty (no cell at all, same
# as any other never-touched cell on the row).
$ws.Range("D4").Clear()
$ws.Range("D7:D18").Clear()

# This sheet becomes the active tab, with D4 selected.
$ws.Activate()
$ws.Range("D4").Select()
